# "Squash: started Landau damping"
# Appends a new "Title and Content" slide at the end of the deck, mirroring
# the layout/formatting of the preceding slide ("Additional steps" - the
# last slide, index 17) so the new slide inherits the same language/locale
# run formatting (en-GB) instead of the engine's generic en-US default.

$p = $ppt.ActivePresentation

# Duplicate the last existing slide so the new slide lands right after it,
# using the same "Title and Content" layout/formatting as the rest of the
# deck.
$lastIndex = $p.Slides.Count
$sourceSlide = $p.Slides.Item($lastIndex)
$dup = $sourceSlide.Duplicate()
$newSlide = $dup.Item(1)

# Replace the title text.
$titleRange = $newSlide.Shapes.Item(1).TextFrame.TextRange
$titleRange.Delete()
$titleRange.Text = "Observation"

# Replace the body text.
$bodyRange = $newSlide.Shapes.Item(2).TextFrame.TextRange
$bodyRange.Delete()
$bodyRange.Text = "I" + [char]0x2019 + "ve been doing 1D3V simulations while in last year it was all 1D1V"
